$d = $word.ActiveDocument

# List of (old, new) text replacements for the summary table figures.
$replacements = @(
    @{ Old = "64 (25.3)";  New = "65 (25.1)" },
    @{ Old = "49 (19.4)";  New = "53 (20.5)" },
    @{ Old = "66 (26.1)";  New = "67 (25.9)" },
    @{ Old = "29 (11.5)";  New = "29 (11.2)" },
    @{ Old = "25 (9.9)";   New = "25 (9.7)" },
    @{ Old = "9 (3.6)";    New = "9 (3.5)" },
    @{ Old = "5 (2.0)";    New = "5 (1.9)" },
    @{ Old = "6 (2.4)";    New = "6 (2.3)" },
    @{ Old = "253 (100.0)"; New = "259 (100.0)" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute(
        $r.Old,   # FindText
        $true,    # MatchCase
        $true,    # MatchWholeWord
        $false,   # MatchWildcards
        $false,   # MatchSoundsLike
        $false,   # MatchAllWordForms
        $true,    # Forward
        1,        # Wrap (wdFindContinue)
        $false,   # Format
        $r.New,   # ReplaceWith
        2         # Replace (wdReplaceAll)
    )
}

$d.Save()
